$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before the old "Data Type" column (column I).
#    This shifts: old I (Data Type) -> J, old J (Reasoning) -> K
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).Insert()

# ---------------------------------------------------------------------------
# 2. Build the two new cell styles we need (bold header / underlined key)
#    by copying the existing bordered style (currently on H2) and tweaking
#    the font - this reuses the existing border definition instead of
#    creating a near duplicate one.
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("I2").Font.Bold = $true        # bold + border -> new style (Table Name header style)

$ws.Range("H2").Copy()
$ws.Range("I3").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("I3").Font.Underline = $true   # underline + border -> new style (key reference style)

# ---------------------------------------------------------------------------
# 3. Header row 2: add the new "Column" heading, re-apply the bold+border
#    style (built above on I2) across the whole header row.
# ---------------------------------------------------------------------------
$ws.Range("I2").Copy()
$ws.Range("H2:K2").PasteSpecial(-4122)   # xlPasteFormats (keep consistent bold+border on the row)

$ws.Range("H2").Value = "Table Name"
$ws.Range("I2").Value = "Column"
$ws.Range("J2").Value = "Data Type"
$ws.Range("K2").Value = "Reasoning"

# ---------------------------------------------------------------------------
# 4. Row 3 (Currency table, 1st column = CurrencyID)
# ---------------------------------------------------------------------------
$ws.Range("I3").Value = "CurrencyID"

# ---------------------------------------------------------------------------
# 5. Rows 4 and 5: new "Column" cells mirror column A values, using the
#    plain bordered style already used by the rest of that row (style of J4/J5).
# ---------------------------------------------------------------------------
$ws.Range("J4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = "CurrencyName"

$ws.Range("J5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = "GDAXEndpoint"

# ---------------------------------------------------------------------------
# 6. New row 7: second copy of the table header (Table Name / Column /
#    Data Type / Reasoning), matching the header already built in row 2.
# ---------------------------------------------------------------------------
$ws.Range("H2:K2").Copy()
$ws.Range("H7:K7").PasteSpecial(-4122)

$ws.Range("H7").Value = "Table Name"
$ws.Range("I7").Value = "Column"
$ws.Range("J7").Value = "Data Type"
$ws.Range("K7").Value = "Reasoning"

# ---------------------------------------------------------------------------
# 7. Row 8 (ExchangeRate table, 1st column = CurrencyID, underlined/key style)
# ---------------------------------------------------------------------------
$ws.Range("I3").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = "CurrencyID"

# ---------------------------------------------------------------------------
# 8. Row 9 (Timestamp) - this one keeps BOTH an underlined "Column" cell
#    (new I9) and the plain one that was shifted from the old I9 into J9.
# ---------------------------------------------------------------------------
$ws.Range("I3").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = "Timestamp"

# ---------------------------------------------------------------------------
# 9. Rows 10-14: remaining "Column" cells mirror column A/C/E/G values with
#    the plain bordered style used across that column.
# ---------------------------------------------------------------------------
$ws.Range("J10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = "Value"

$ws.Range("J11").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = "Growth"

$ws.Range("J12").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = "GOFAIPredictedGrowth"

$ws.Range("J13").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = "NeuralNetworkPredictedGrowth"

$ws.Range("J14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = "LastGDAXTrade"

# ---------------------------------------------------------------------------
# 10. New explanatory notes below the tables (rows 16-18, plain, no style).
# ---------------------------------------------------------------------------
$ws.Range("H16").Value = "A Currency has a single Exchange Rate at any given point in time."
$ws.Range("H17").Value = "The Exchange Rate's time specifies value (in USD for simplicity), growth and predictions."
$ws.Range("H18").Value = "The LastGDAXTrade number is stored for ease of navigation to gaps in data should they occur."

# ---------------------------------------------------------------------------
# 11. Re-fit the new column (matching the width Excel's own "best fit"
#     produced for the final layout) and move the active selection,
#     matching the workbook's final on-screen state.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 29.6

$ws.Range("G22").Select()
